$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13..130 down to 14..131
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new record's data
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13,3).Value = "Bíobío"
$ws.Cells.Item(13,4).Value = 44901
$ws.Cells.Item(13,5).Value = 8
$ws.Cells.Item(13,6).Value = 100112001
$ws.Cells.Item(13,7).Value = "Berenjena"
$ws.Cells.Item(13,8).Value = "Sin especificar"
$ws.Cells.Item(13,9).Value = "Primera"
$ws.Cells.Item(13,10).Value = 220
$ws.Cells.Item(13,11).Value = 11000
$ws.Cells.Item(13,12).Value = 12000
$ws.Cells.Item(13,13).Value = 11455
$ws.Cells.Item(13,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(13,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13,16).Value = 191
$ws.Cells.Item(13,17).Value = 60
$ws.Cells.Item(13,18).Value = "Hortaliza"
